# Slide 8 ("Git basics"): the bullet list text box ("Объект 2") paragraph
#   "$ git push – Push all local commits to remote repository"
# gets its last run (" push – Push all local commits to remote repository")
# split into three runs so the line reads:
#   "$ git push origin (+)HEAD – Push all local commits to remote repository"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange

# Paragraph 9 is "$ git push – Push all local commits to remote repository".
$para = $tr.Paragraphs(9, 1)

# The run we need to edit is " push – Push all local commits to remote
# repository" (51 characters), right after "$ git". Leave its leading space
# (local position 6, i.e. "$ git" + 1) untouched as its own run, and only
# rewrite from "push" onward (local position 7), which is 50 characters.
$wordStart = $para.Start + 6
$wordLen = 50

$tail = $tr.Characters($wordStart, $wordLen)
$tail.Text = "push origin (+)HEAD " + [char]0x2013 + " Push all local commits to remote repository"

# Split the freshly written text into two runs: "push origin (+)HEAD " and
# "– Push all local commits to remote repository".
$head = $tr.Characters($wordStart, 20)
$head.Text = "push origin (+)HEAD "
